$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) cells with refreshed market data.
# Numeric-looking price strings get a leading apostrophe so Excel keeps
# them as literal text (matching the original "Price" column formatting)
# instead of silently converting them to floating point numbers.
$ws.Range('D2').Value = '62.677.43'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '3.195.65'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''593.69'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').Value = '''135.98'
$ws.Range('E6').Value = '  -5.50%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.195.29'
$ws.Range('E8').Value = '  -3.36%  '
$ws.Range('D9').Value = '''0.506'
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').Value = '''0.455'
$ws.Range('E12').Value = '  -4.26%  '
$ws.Range('D13').Value = '''0.0000238'
$ws.Range('E13').Value = '  -4.43%  '
$ws.Range('D14').Value = '''33.52'
$ws.Range('E14').Value = '  -4.28%  '
$ws.Range('D15').Value = '3.724.38'
$ws.Range('E15').Value = '  -3.44%  '
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '3.198.32'
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').Value = '62.671.83'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('D19').Value = '''6.70'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').Value = '''463.60'
$ws.Range('E20').Value = '  -4.35%  '
$ws.Range('D21').Value = '''13.90'
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('D22').Value = '''0.714'
$ws.Range('E22').Value = '  -4.34%  '
$ws.Range('D23').Value = '''7.67'
$ws.Range('E23').Value = '  -5.17%  '
$ws.Range('D25').Value = '''83.77'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').Value = '''7.87'
$ws.Range('E29').Value = '  -5.33%  '
$ws.Range('D30').Value = '''6.90'
$ws.Range('E30').Value = '  -5.64%  '
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('D32').Value = '''27.54'
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('E33').Value = '  -4.85%  '
$ws.Range('D34').Value = '''2.43'
$ws.Range('D35').Value = '''1.05'
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('D36').Value = '''5.85'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('D37').Value = '''51.52'
$ws.Range('E37').Value = '  -3.60%  '
$ws.Range('D38').Value = '0.0₃0699'
$ws.Range('E38').Value = '  -5.51%  '
$ws.Range('D39').Value = '''0.0389'
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('D40').Value = '''419.61'
$ws.Range('E40').Value = '  -2.97%  '
$ws.Range('D41').Value = '2.998.18'
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('D43').Value = '''8.10'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('D44').Value = '''2.62'
$ws.Range('E44').Value = '  -6.18%  '
$ws.Range('D45').Value = '''0.254'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('E46').Value = '  -5.05%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''125.62'
$ws.Range('E48').Value = '  +1.63%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = '''35.32'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''25.64'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('E51').Value = '  -3.13%  '
